$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.000002317355952907718
$ws.Range("C2").Value = 0.04240448674262143
$ws.Range("D2").Value = 3.900430680208489
$ws.Range("E2").Value = 645.3272768299601
$ws.Range("G2").Value = 649.2701143142672
